$d = $word.ActiveDocument

# 1) Insert the missing space between the closing "+++" of the
#    tratamiento/E ternary template and "del objeto y resultado..."
#    enterada`+++del objeto  ->  enterada`+++ del objeto
$d.Content.Find.Execute(
    'enterada`+++del objeto',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'enterada`+++ del objeto',
    2)

# 2) Turn the plain word "esta" (right after "Senor Notario, ") into the
#    templated gender ternary expression, same pattern used elsewhere in
#    the document:
#    Senor Notario, esta se afirma
#      ->
#    Senor Notario, +++= documento.otorgantes[0].tratamiento >= `E`? `este`:`esta`+++ se afirma
$d.Content.Find.Execute(
    'Notario, esta se afirma',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'Notario, +++= documento.otorgantes[0].tratamiento >= `E`? `este`:`esta`+++ se afirma',
    2)
